$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the period-header row (row 1, columns B:F) to the new
# lowercase/underscored labels.
$ws.Range("B1").Value = "set_2024_a_ago_2025"
$ws.Range("C1").Value = "mai_2025_a_abr_2025"
$ws.Range("D1").Value = "jan_2024_a_dez_2024"
$ws.Range("E1").Value = "set_2023_a_ago_2024"
$ws.Range("F1").Value = "mai_2023_a_abr_2024"

# Move the active selection to F1, matching the saved view state.
$ws.Range("F1").Select()
